$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect the Price column (D) so numeric-looking text (e.g. '11.38') is
# written back as text, matching the source inlineStr cells, instead of
# being auto-coerced to a number by Excel's smart Value assignment.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range("D2").Value = "68.313.39"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "2.646.83"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "597.41"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").Value = "155.92"
$ws.Range("E6").Value = "  +0.94%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("E9").Value = "  +6.74%  "
$ws.Range("D10").Value = "0.157"
$ws.Range("E10").Value = "  -0.94%  "
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("D13").Value = "28.08"
$ws.Range("E13").Value = "  +1.67%  "
$ws.Range("E14").Value = "  +2.11%  "
$ws.Range("D15").Value = "3.126.58"
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("D16").Value = "68.218.39"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").Value = "2.643.70"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").Value = "11.38"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "363.90"
$ws.Range("E19").Value = "  -1.60%  "
$ws.Range("D20").Value = "7.47"
$ws.Range("E20").Value = "  +1.07%  "
$ws.Range("E21").Value = "  +3.33%  "
$ws.Range("D22").Value = "4.84"
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("E23").Value = "  -0.69%  "
$ws.Range("E24").Value = "  +3.17%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  -1.78%  "
$ws.Range("E27").Value = "  +1.57%  "
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("D30").Value = "559.63"
$ws.Range("E30").Value = "  -2.15%  "
$ws.Range("D31").Value = "8.04"
$ws.Range("E31").Value = "  +1.02%  "
$ws.Range("E32").Value = "  +1.24%  "
$ws.Range("E33").Value = "  +0.33%  "
$ws.Range("E34").Value = "  +1.66%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  +3.41%  "
$ws.Range("D37").Value = "161.29"
$ws.Range("E37").Value = "  +0.60%  "
$ws.Range("D38").Value = "19.33"
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("E39").Value = "  +1.50%  "
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("D41").Value = "5.34"
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("D42").Value = "0.0₆0340"
$ws.Range("E42").Value = "  +4.51%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "17.77"
$ws.Range("E44").Value = "  +0.91%  "
$ws.Range("D46").Value = "40.35"
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("D47").Value = "159.43"
$ws.Range("E47").Value = "  +2.83%  "
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("D49").Value = "21.97"
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.0785"
$ws.Range("E50").Value = "  +0.96%  "
$ws.Range("B51").Value = "Optimism"
$ws.Range("C51").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D51").Value = "1.69"
$ws.Range("E51").Value = "  -0.23%  "

# Remove the temporary text-format override so styling matches the original
# (no explicit number format on these cells).
$priceCol.ClearFormats()
